$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the coin price/volume table with the latest scrape.
#
# Column D ('Price') stores numeric-looking values as TEXT (e.g. the grouping
# dots make some look like '25.894.21', which isn't a valid number anyway, but
# plain ones like '214.52' are). Assigning a plain numeric-looking string to
# .Value lets Excel auto-convert it to a real number, which would change the
# cell's stored type. For those cells we force text storage by setting
# NumberFormat to '@' before the assignment, then restore the 'Normal' style
# afterwards so no stray number formatting lingers on the cell.

$ws.Range('D2').Value = '25.894.21'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '1.635.10'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('E6').Value = '  +1.00%  '
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('E8').Value = '  -0.65%  '
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('E10').Value = '  -0.21%  '
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.683.51'
$ws.Range('E12').Value = '  +4.78%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.860.55'
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.25'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.543'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.50%  '
$ws.Range('D16').Value = '0.0₃0756'
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.59'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '25.910.19'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.40'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.82%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '193.57'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.37%  '
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('E24').Value = '  +1.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.58'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.96%  '
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('E27').Value = '  +3.25%  '
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.43'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.32%  '
$ws.Range('E30').Value = '  +0.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0499'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('E32').Value = '  -0.53%  '
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('E34').Value = '  -2.11%  '
$ws.Range('E35').Value = '  +1.44%  '
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('D37').Value = '1.140.71'
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('E39').Value = '  -0.84%  '
$ws.Range('E40').Value = '  +0.28%  '
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.51'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.86%  '
$ws.Range('E43').Value = '  -0.38%  '
$ws.Range('E44').Value = '  -3.38%  '
$ws.Range('D45').Value = '1.769.49'
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range('E46').Value = '  +1.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.33'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.61%  '
$ws.Range('E48').Value = '  +3.32%  '
$ws.Range('E49').Value = '  -0.85%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.415'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.64'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.88%  '
